$d = $word.ActiveDocument

# Helper: force a run-split boundary at both ends of the given range by
# toggling Bold off->on (or on->off) and back to its original value. The
# underlying engine re-renders run boundaries for any range whose direct
# formatting changes, so flipping a property and restoring it is enough to
# carve the range out into its own run(s) without altering its appearance.
function Split-At($pos1, $pos2) {
    $rg = $d.Range($pos1, $pos2)
    $orig = $rg.Bold
    if ($orig -eq 1 -or $orig -eq -1) {
        $rg.Bold = 0
        $rg.Bold = $orig
    } else {
        $rg.Bold = 1
        $rg.Bold = 0
    }
}

# ---------------------------------------------------------------------
# Edit 1: "Moving the ball: from left and right, up and down"
#      -> "Moving the ball: from left and right, " / "up," / " and down"
#         (the trailing " until it bounces to the wall" run is untouched)
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Moving the ball: from left and right, up and down", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'Moving the ball' sentence"
}

$full1 = "Moving the ball: from left and right, up and down"
$idx1 = $full1.IndexOf("up and down")
$R1 = $r1.Start

$upStart = $R1 + $idx1
$upEnd = $upStart + 2
$commaSub = $d.Range($upStart, $upEnd)
$commaSub.InsertAfter(",")

# Boundaries for the resulting runs (absolute positions, after the comma insert):
#   run1: "Moving the ball: from left and right, "   [R1, upStart)
#   run2: "up,"                                       [upStart, upEnd+1)
#   run3: " and down"                                 [upEnd+1, R1+full1.Length+1)
$q1 = $R1
$q2 = $upStart
$q3 = $upEnd + 1
$q4 = $R1 + $full1.Length + 1

Split-At $q1 $q2
Split-At $q2 $q3
Split-At $q3 $q4

# ---------------------------------------------------------------------
# Edit 2: "The game will end once you click the quit button"
#      -> "The game will end once you click the quit " / "button,"
#         (the leading spaces run and the trailing " or the time has
#         expired" / "." runs are untouched)
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("The game will end once you click the quit button", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'The game will end' sentence"
}

$full2 = "The game will end once you click the quit button"
$idx2 = $full2.IndexOf("button")
$R2 = $r2.Start

$buttonStart = $R2 + $idx2
$buttonEnd = $R2 + $full2.Length
$commaSub2 = $d.Range($buttonStart, $buttonEnd)
$commaSub2.InsertAfter(",")

# Boundaries for the resulting runs (absolute positions, after the comma insert):
#   run0: "       " (7 leading spaces)                [R2-7, R2)
#   run1: "The game will end once you click the quit " [R2, buttonStart)
#   run2: "button,"                                    [buttonStart, buttonEnd+1)
#   run3: " or the time has expired"                   [buttonEnd+1, buttonEnd+1+25)
#   run4: "."                                          [buttonEnd+1+25, buttonEnd+1+25+1)
$p0 = $R2 - 7
$p1 = $R2
$p2 = $buttonStart
$p3 = $buttonEnd + 1
$orLen = " or the time has expired".Length
$p4 = $p3 + $orLen
$p5 = $p4 + 1

Split-At $p0 $p1
Split-At $p1 $p2
Split-At $p2 $p3
Split-At $p3 $p4
Split-At $p4 $p5

Write-Output "Done"
